# Auto-generated script applying scheduled price-runner updates to the
# FFXIV leve profit workbook (columns H-N: currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 152
$ws.Range("I2").Value = 152
$ws.Range("K2").Value = 152
$ws.Range("M2").Value = -39
$ws.Range("H17").Value = 1590
$ws.Range("J17").Value = 1590
$ws.Range("L17").Value = 4770
$ws.Range("N17").Value = -5106
$ws.Range("H43").Value = 3402
$ws.Range("J43").Value = 3402
$ws.Range("L43").Value = 3402
$ws.Range("N43").Value = -3540
$ws.Range("H64").Value = 3995
$ws.Range("I64").Value = 3995
$ws.Range("K64").Value = 3995
$ws.Range("M64").Value = -3747
$ws.Range("H67").Value = 3995
$ws.Range("I67").Value = 3995
$ws.Range("K67").Value = 3995
$ws.Range("M67").Value = -3137
$ws.Range("H74").Value = 3716.6667
$ws.Range("I74").Value = 3660
$ws.Range("K74").Value = 3660
$ws.Range("M74").Value = -2724
$ws.Range("H77").Value = 3716.6667
$ws.Range("I77").Value = 3660
$ws.Range("K77").Value = 18300
$ws.Range("M77").Value = -13620
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H111").Value = 1997.5
$ws.Range("J111").Value = 1995
$ws.Range("L111").Value = 5985
$ws.Range("N111").Value = -12119
$ws.Range("H127").Value = 448.75
$ws.Range("I127").Value = 448.75
$ws.Range("K127").Value = 1346.25
$ws.Range("M127").Value = 3613.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61248
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186240
$ws.Range("H74").Value = 2999.8333
$ws.Range("I74").Value = 2099.8
$ws.Range("K74").Value = 2099.8
$ws.Range("M74").Value = -1225.8
$ws.Range("H77").Value = 2999.8333
$ws.Range("I77").Value = 2099.8
$ws.Range("K77").Value = 10499
$ws.Range("M77").Value = -6131
$ws.Range("H132").Value = 6200.6
$ws.Range("I132").Value = 3001
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 9003
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -6473
$ws.Range("N132").Value = -38060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6505.5
$ws.Range("I134").Value = 1011
$ws.Range("K134").Value = 3033
$ws.Range("M134").Value = -498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 435.4
$ws.Range("I5").Value = 392.42856
$ws.Range("J5").Value = 535.6667
$ws.Range("K5").Value = 392.42856
$ws.Range("L5").Value = 535.6667
$ws.Range("M5").Value = -280.42856
$ws.Range("N5").Value = -759.6667
$ws.Range("H132").Value = 4038.1428
$ws.Range("I132").Value = 1253.6
$ws.Range("K132").Value = 3760.8
$ws.Range("M132").Value = -1230.8
$ws.Range("H134").Value = 6999.5
$ws.Range("I134").Value = 6999.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 20998.5
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -18463.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7889166.5
$ws.Range("I11").Value = 11750000
$ws.Range("J11").Value = 4800500
$ws.Range("K11").Value = 11750000
$ws.Range("L11").Value = 4800500
$ws.Range("M11").Value = -11749861
$ws.Range("N11").Value = -4800778
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H102").Value = 360.08334
$ws.Range("I102").Value = 272.2
$ws.Range("J102").Value = 799.5
$ws.Range("K102").Value = 272.2
$ws.Range("L102").Value = 799.5
$ws.Range("M102").Value = 1349.8
$ws.Range("N102").Value = -4043.5
$ws.Range("H132").Value = 2642.2856
$ws.Range("I132").Value = 1999.4
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 5998.200000000001
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -3468.200000000001
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 576.1667
$ws.Range("I22").Value = 551
$ws.Range("J22").Value = 702
$ws.Range("K22").Value = 551
$ws.Range("L22").Value = 702
$ws.Range("M22").Value = -256
$ws.Range("N22").Value = -1292
$ws.Range("H27").Value = 576.1667
$ws.Range("I27").Value = 551
$ws.Range("J27").Value = 702
$ws.Range("K27").Value = 551
$ws.Range("L27").Value = 702
$ws.Range("M27").Value = -444
$ws.Range("N27").Value = -916
$ws.Range("H132").Value = 4448.7
$ws.Range("I132").Value = 3784
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 11352
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -8822
$ws.Range("N132").Value = -23058.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28297.75
$ws.Range("I45").Value = 22570
$ws.Range("K45").Value = 22570
$ws.Range("M45").Value = -22079
$ws.Range("H62").Value = 5555
$ws.Range("J62").Value = 5555
$ws.Range("L62").Value = 5555
$ws.Range("N62").Value = -6803
$ws.Range("H65").Value = 5555
$ws.Range("J65").Value = 5555
$ws.Range("L65").Value = 27775
$ws.Range("N65").Value = -34015
$ws.Range("H81").Value = 966.6667
$ws.Range("I81").Value = 850
$ws.Range("K81").Value = 1700
$ws.Range("M81").Value = -639
$ws.Range("H84").Value = 966.6667
$ws.Range("I84").Value = 850
$ws.Range("K84").Value = 8500
$ws.Range("M84").Value = -3196
$ws.Range("H136").Value = 1514.4117
$ws.Range("I136").Value = 1654.5
$ws.Range("K136").Value = 4963.5
$ws.Range("M136").Value = -2413.5
